# Apply scraped-schedule update for Linea 141 (commit: "Horarios actualizados Linea 141 - 1279")
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Sheet 1 ---
$ws1.Cells.Item(2,1).Value = "Última actualización: 08:49:06"
$ws1.Cells.Item(3,1).Value = "Total filas: 83"
$ws1.Cells.Item(23,1).Value = "06:15:23"
$ws1.Cells.Item(23,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(23,4).Value = 66
$ws1.Cells.Item(24,1).Value = "06:46:40"
$ws1.Cells.Item(24,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(24,4).Value = 35
$ws1.Cells.Item(45,1).Value = "08:14:55"
$ws1.Cells.Item(45,3).Value = "15_ABASTO"
$ws1.Cells.Item(45,4).Value = 15
$ws1.Cells.Item(46,1).Value = "07:51:40"
$ws1.Cells.Item(46,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(46,4).Value = 38
$ws1.Cells.Item(52,1).Value = "08:49:06"
$ws1.Cells.Item(52,2).Value = "08:50"
$ws1.Cells.Item(52,3).Value = "215C_EL PATO"
$ws1.Cells.Item(52,4).Value = 1
$ws1.Cells.Item(53,1).Value = "08:14:55"
$ws1.Cells.Item(53,2).Value = "08:51"
$ws1.Cells.Item(53,4).Value = 37
$ws1.Cells.Item(54,1).Value = "06:58:58"
$ws1.Cells.Item(54,2).Value = "08:52"
$ws1.Cells.Item(54,4).Value = 114
$ws1.Cells.Item(56,1).Value = "08:49:06"
$ws1.Cells.Item(56,2).Value = "08:53"
$ws1.Cells.Item(56,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(56,4).Value = 4
$ws1.Cells.Item(57,1).Value = "08:49:06"
$ws1.Cells.Item(57,2).Value = "08:54"
$ws1.Cells.Item(57,3).Value = "215B_EL PATO"
$ws1.Cells.Item(57,4).Value = 5
$ws1.Cells.Item(58,1).Value = "08:14:55"
$ws1.Cells.Item(58,2).Value = "08:57"
$ws1.Cells.Item(58,4).Value = 43
$ws1.Cells.Item(59,1).Value = "08:49:06"
$ws1.Cells.Item(59,2).Value = "08:58"
$ws1.Cells.Item(59,3).Value = "215A_EL PATO"
$ws1.Cells.Item(59,4).Value = 9
$ws1.Cells.Item(60,1).Value = "08:14:55"
$ws1.Cells.Item(60,2).Value = "09:04"
$ws1.Cells.Item(60,4).Value = 50
$ws1.Cells.Item(61,1).Value = "08:49:06"
$ws1.Cells.Item(61,2).Value = "09:05"
$ws1.Cells.Item(61,3).Value = "10_OLMOS"
$ws1.Cells.Item(61,4).Value = 16
$ws1.Cells.Item(62,1).Value = "08:49:06"
$ws1.Cells.Item(62,2).Value = "09:06"
$ws1.Cells.Item(62,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(62,4).Value = 17
$ws1.Cells.Item(63,1).Value = "07:26:49"
$ws1.Cells.Item(63,2).Value = "09:16"
$ws1.Cells.Item(63,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(63,4).Value = 110
$ws1.Cells.Item(64,1).Value = "08:14:55"
$ws1.Cells.Item(64,3).Value = "14_ABASTO"
$ws1.Cells.Item(64,4).Value = 63
$ws1.Cells.Item(65,1).Value = "08:49:06"
$ws1.Cells.Item(65,2).Value = "09:17"
$ws1.Cells.Item(65,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(65,4).Value = 28
$ws1.Cells.Item(66,1).Value = "08:49:06"
$ws1.Cells.Item(66,4).Value = 29
$ws1.Cells.Item(67,1).Value = "08:49:06"
$ws1.Cells.Item(67,2).Value = "09:18"
$ws1.Cells.Item(67,3).Value = "15X38_ABASTO"
$ws1.Cells.Item(67,4).Value = 29
$ws1.Cells.Item(68,1).Value = "07:51:40"
$ws1.Cells.Item(68,2).Value = "09:21"
$ws1.Cells.Item(68,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(68,4).Value = 90
$ws1.Cells.Item(69,1).Value = "08:14:55"
$ws1.Cells.Item(69,2).Value = "09:28"
$ws1.Cells.Item(69,4).Value = 74
$ws1.Cells.Item(70,1).Value = "08:49:06"
$ws1.Cells.Item(70,2).Value = "09:29"
$ws1.Cells.Item(70,3).Value = "10_OLMOS"
$ws1.Cells.Item(70,4).Value = 40
$ws1.Cells.Item(71,1).Value = "08:49:06"
$ws1.Cells.Item(71,2).Value = "09:31"
$ws1.Cells.Item(71,4).Value = 42
$ws1.Cells.Item(72,2).Value = "09:31"
$ws1.Cells.Item(72,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(72,4).Value = 77
$ws1.Cells.Item(73,2).Value = "09:33"
$ws1.Cells.Item(73,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(73,4).Value = 58
$ws1.Cells.Item(74,2).Value = "09:35"
$ws1.Cells.Item(74,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(74,4).Value = 81
$ws1.Cells.Item(75,1).Value = "08:49:06"
$ws1.Cells.Item(75,2).Value = "09:39"
$ws1.Cells.Item(75,3).Value = "15_ABASTO"
$ws1.Cells.Item(75,4).Value = 50
$ws1.Cells.Item(76,1).Value = "08:14:55"
$ws1.Cells.Item(76,2).Value = "09:41"
$ws1.Cells.Item(76,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(76,4).Value = 87
$ws1.Cells.Item(77,1).Value = "08:49:06"
$ws1.Cells.Item(77,2).Value = "09:42"
$ws1.Cells.Item(77,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(77,4).Value = 53
$ws1.Cells.Item(78,1).Value = "08:49:06"
$ws1.Cells.Item(78,2).Value = "09:43"
$ws1.Cells.Item(78,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(78,4).Value = 54
$ws1.Cells.Item(79,1).Value = "08:49:06"
$ws1.Cells.Item(79,2).Value = "09:53"
$ws1.Cells.Item(79,3).Value = "10_OLMOS"
$ws1.Cells.Item(79,4).Value = 64
$ws1.Cells.Item(80,2).Value = "09:58"
$ws1.Cells.Item(80,3).Value = "215C_EL PATO"
$ws1.Cells.Item(80,4).Value = 104
$ws1.Cells.Item(81,1).Value = "08:49:06"
$ws1.Cells.Item(81,2).Value = "09:59"
$ws1.Cells.Item(81,3).Value = "215C_EL PATO"
$ws1.Cells.Item(81,4).Value = 70
$ws1.Cells.Item(82,1).Value = "08:14:55"
$ws1.Cells.Item(82,2).Value = "10:05"
$ws1.Cells.Item(82,3).Value = "14_ABASTO"
$ws1.Cells.Item(82,4).Value = 111
$ws1.Cells.Item(83,1).Value = "08:49:06"
$ws1.Cells.Item(83,2).Value = "10:06"
$ws1.Cells.Item(83,3).Value = "14_ABASTO"
$ws1.Cells.Item(83,4).Value = 77
$ws1.Cells.Item(84,1).Value = "08:49:06"
$ws1.Cells.Item(84,2).Value = "10:13"
$ws1.Cells.Item(84,3).Value = "17X38_ROMERO"
$ws1.Cells.Item(84,4).Value = 84
$ws1.Cells.Item(85,1).Value = "08:49:06"
$ws1.Cells.Item(85,2).Value = "10:25"
$ws1.Cells.Item(85,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(85,4).Value = 96
$ws1.Cells.Item(85,5).Value = "LP1912"
$ws1.Cells.Item(86,1).Value = "08:49:06"
$ws1.Cells.Item(86,2).Value = "10:29"
$ws1.Cells.Item(86,3).Value = "15_ABASTO"
$ws1.Cells.Item(86,4).Value = 100
$ws1.Cells.Item(86,5).Value = "LP1912"
$ws1.Cells.Item(87,1).Value = "08:49:06"
$ws1.Cells.Item(87,2).Value = "10:44"
$ws1.Cells.Item(87,3).Value = "11X44_ETCHEVERRY"
$ws1.Cells.Item(87,4).Value = 115
$ws1.Cells.Item(87,5).Value = "LP1912"
$ws1.Cells.Item(88,1).Value = "08:49:06"
$ws1.Cells.Item(88,2).Value = "10:46"
$ws1.Cells.Item(88,3).Value = "15_P INDUSTRIAL"
$ws1.Cells.Item(88,4).Value = 117
$ws1.Cells.Item(88,5).Value = "LP1912"

# --- Sheet 2 ---
$ws2.Cells.Item(2,1).Value = "Última actualización: 08:49:06"
$ws2.Cells.Item(3,1).Value = "Total filas: 16"
$ws2.Cells.Item(15,1).Value = "08:49:06"
$ws2.Cells.Item(15,2).Value = "08:50"
$ws2.Cells.Item(15,3).Value = "215C_EL PATO"
$ws2.Cells.Item(15,4).Value = 1
$ws2.Cells.Item(16,1).Value = "08:14:55"
$ws2.Cells.Item(16,2).Value = "08:53"
$ws2.Cells.Item(16,4).Value = 39
$ws2.Cells.Item(17,1).Value = "08:49:06"
$ws2.Cells.Item(17,2).Value = "08:54"
$ws2.Cells.Item(17,3).Value = "215B_EL PATO"
$ws2.Cells.Item(17,4).Value = 5
$ws2.Cells.Item(18,1).Value = "08:14:55"
$ws2.Cells.Item(18,2).Value = "08:57"
$ws2.Cells.Item(18,4).Value = 43
$ws2.Cells.Item(19,1).Value = "08:49:06"
$ws2.Cells.Item(19,2).Value = "08:58"
$ws2.Cells.Item(19,3).Value = "215A_EL PATO"
$ws2.Cells.Item(19,4).Value = 9
$ws2.Cells.Item(20,1).Value = "08:14:55"
$ws2.Cells.Item(20,2).Value = "09:58"
$ws2.Cells.Item(20,4).Value = 104
$ws2.Cells.Item(21,1).Value = "08:49:06"
$ws2.Cells.Item(21,2).Value = "09:59"
$ws2.Cells.Item(21,3).Value = "215C_EL PATO"
$ws2.Cells.Item(21,4).Value = 70
$ws2.Cells.Item(21,5).Value = "LP1912"

# --- Sheet 3 ---
$ws3.Cells.Item(2,1).Value = "Última actualización: 08:49:06"
$ws3.Cells.Item(3,1).Value = "Total filas: 13"
$ws3.Cells.Item(12,1).Value = "08:49:06"
$ws3.Cells.Item(12,2).Value = "08:52"
$ws3.Cells.Item(12,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(12,4).Value = 3
$ws3.Cells.Item(12,5).Value = "L6203"
$ws3.Cells.Item(13,1).Value = "08:14:55"
$ws3.Cells.Item(13,2).Value = "09:20"
$ws3.Cells.Item(13,4).Value = 66
$ws3.Cells.Item(14,1).Value = "08:49:06"
$ws3.Cells.Item(14,2).Value = "09:21"
$ws3.Cells.Item(14,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(14,4).Value = 32
$ws3.Cells.Item(14,5).Value = "L6173"
$ws3.Cells.Item(15,1).Value = "08:14:55"
$ws3.Cells.Item(15,2).Value = "10:12"
$ws3.Cells.Item(15,4).Value = 118
$ws3.Cells.Item(16,1).Value = "08:49:06"
$ws3.Cells.Item(16,2).Value = "10:13"
$ws3.Cells.Item(16,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(16,4).Value = 84
$ws3.Cells.Item(16,5).Value = "L6203"
$ws3.Cells.Item(17,1).Value = "08:49:06"
$ws3.Cells.Item(17,2).Value = "10:30"
$ws3.Cells.Item(17,3).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(17,4).Value = 101
$ws3.Cells.Item(17,5).Value = "L6173"
$ws3.Cells.Item(18,1).Value = "08:49:06"
$ws3.Cells.Item(18,2).Value = "10:31"
$ws3.Cells.Item(18,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(18,4).Value = 102
$ws3.Cells.Item(18,5).Value = "L6173"
